$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values (repetition indices)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values
$ws.Range("B2").Value = 60.780713310805943
$ws.Range("C2").Value = 49.846522874245437
$ws.Range("D2").Value = 65.197813094466525
$ws.Range("E2").Value = 52.598329430698286

# Row 3 data values
$ws.Range("B3").Value = 63.221206623705854
$ws.Range("C3").Value = 44.659973050356776
$ws.Range("D3").Value = 74.969337591465788
$ws.Range("E3").Value = 50.060005343044153

# Update the selection to match the new narrower highlighted range
$ws.Range("B1:E3").Select()
